{"js": "// Update the meeting date from 3/1 to 3/8 in the first paragraph that\n// contains the \"Meeting time/place\" text, then move the \"_GoBack\"\n// bookmark (originally on the last, empty paragraph of the document)\n// to sit right after that updated run, matching what Word itself does\n// when a user last edits text in that spot.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The meeting info paragraph is the second paragraph in the document\n// (index 1) and starts with \"Meeting time/place: 3/1...\".\nconst meetingParagraph = paragraphs.items[1];\n\nconst found = meetingParagraph.search(\"Meeting time/place: 3/1\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nconst target = found.items[0];\ntarget.insertText(\"Meeting time/place: 3/8\", \"Replace\");\nawait context.sync();\n\n// Remove the pre-existing \"_GoBack\" bookmark (it lived on the trailing\n// empty paragraph) before re-adding it in its new location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-insert \"_GoBack\" immediately after the text we just replaced.\nconst afterEdit = target.getRange(\"End\");\nafterEdit.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Update the meeting date from 3/1 to 3/8 in the \"Meeting time/place\"\n# line, then move the \"_GoBack\" bookmark (previously sitting on the\n# trailing empty paragraph at the end of the document) so it lands\n# right after the updated text -- matching where Word itself leaves\n# that bookmark after the last text edit.\n\n$d = $word.ActiveDocument\n\n# 1. Locate \"Meeting time/place: 3/1\" precisely.\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \"Meeting time/place: 3/1\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.MatchCase = $true\n$find.Execute() | Out-Null\n\n$endOfMatch = $searchRange.End\n\n# 2. Relocate the \"_GoBack\" bookmark to sit right after the matched\n#    text BEFORE editing anything, so the text edit below only touches\n#    the run text itself and does not merge across the bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$bmRange = $d.Range($endOfMatch, $endOfMatch)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# 3. Replace just the trailing \"1\" in \"3/1\" with \"8\" so it reads \"3/8\".\n$lastCharRange = $d.Range($endOfMatch - 1, $endOfMatch)\n$lastCharRange.Text = \"8\"\n"}
